$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.204.74"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.426.39"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.74"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.94"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.431.07"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  -8.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.424"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.021.71"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.12"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("E16").Value = "  -6.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.272.83"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.424.37"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.60"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.50"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.55"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.517"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.04"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.98"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.08"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.43"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  +11.59%  "
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.812.03"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0729"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0304"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.91"
$ws.Range("E47").Value = "  +8.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +7.67%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "43.08"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.84"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.48"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.31"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.32"
$ws.Range("E51").Value = "  -3.04%  "
